# Add data for 2022-03-25 (carjacking by month, YoY historical)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet & update the "as-of" header label
$ws.Name = "Through 2022-03-17"
$ws.Range("I1").Value = "2022 (through 03-17)"

# Update the 2022 monthly figures (column I) and the Total row (row 14)
$ws.Range("I4").Value = 74
$ws.Range("H8").Value = 150
$ws.Range("H14").Value = 1852
$ws.Range("I14").Value = 374
